$d = $word.ActiveDocument

# --- Text edits in clause 4.3 -------------------------------------------
# Run 1: insert "6 tiếng so với" before "thời gian dự kiến" and append
# "% giá gói thuê" after "${penalty}". (Single-quoted so PowerShell does
# not try to interpolate the "${penalty}" token.)
$search1 = '4.3 Nếu bên B trả xe lố thời gian dự kiến sẽ phải chịu thêm ${penalty}'
$replace1 = '4.3 Nếu bên B trả xe lố 6 tiếng so với thời gian dự kiến sẽ phải chịu thêm ${penalty}% giá gói thuê'
$d.Content.Find.Execute(
    $search1,
    $true, $false, $false, $false, $false, $true, 1, $false,
    $replace1,
    2) | Out-Null

# Run 2: the trailing "% giá gói thuê đã chọn." collapses to just "."
$search2 = '% giá gói thuê đã chọn.'
$replace2 = '.'
$d.Content.Find.Execute(
    $search2,
    $true, $false, $false, $false, $false, $true, 1, $false,
    $replace2,
    2) | Out-Null

# --- Style tweaks: mark Table Grid / Light Shading Accent 1 / Light
#     Shading Accent 6 as QuickStyles (adds <w:qFormat/> both on the
#     <w:style> element and its latent-style exception entry). ------------
$d.Styles.Item("Table Grid").QuickStyle = $true
$d.Styles.Item("Light Shading Accent 1").QuickStyle = $true
$d.Styles.Item("Light Shading Accent 6").QuickStyle = $true

# --- Drop the redundant cell-margin override on the first row of the
#     first table (tblPrEx/tblCellMar); tblBorders stays. ----------------
$t = $d.Tables.Item(1)
$r = $t.Rows.Item(1)
